$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 4390901
$ws.Range("I86").Value = 4333.75
$ws.Range("J86").Value = 6584185
$ws.Range("K86").Value = 4333.75
$ws.Range("L86").Value = 6584185
$ws.Range("M86").Value = -3210.75
$ws.Range("N86").Value = -6586431
# Row 89
$ws.Range("H89").Value = 4390901
$ws.Range("I89").Value = 4333.75
$ws.Range("J89").Value = 6584185
$ws.Range("K89").Value = 21668.75
$ws.Range("L89").Value = 32920925
$ws.Range("M89").Value = -16052.75
$ws.Range("N89").Value = -32932157
# Row 113
$ws.Range("H113").Value = 15686.875
$ws.Range("I113").Value = 13100.8
$ws.Range("J113").Value = 19997
$ws.Range("K113").Value = 13100.8
$ws.Range("L113").Value = 19997
$ws.Range("M113").Value = -9846.799999999999
$ws.Range("N113").Value = -26505
# Row 116
$ws.Range("H116").Value = 6790.28
$ws.Range("I116").Value = 7018.0713
$ws.Range("J116").Value = 6500.364
$ws.Range("K116").Value = 7018.0713
$ws.Range("L116").Value = 6500.364
$ws.Range("M116").Value = -3576.0713
$ws.Range("N116").Value = -13384.364
# Row 137
$ws.Range("H137").Value = 2718.12
$ws.Range("I137").Value = 1953.4445
$ws.Range("K137").Value = 5860.333500000001
$ws.Range("M137").Value = -3310.333500000001
# Row 138
$ws.Range("H138").Value = 6355.265
$ws.Range("I138").Value = 3640.9333
$ws.Range("J138").Value = 7552.7646
$ws.Range("K138").Value = 10922.7999
$ws.Range("L138").Value = 22658.2938
$ws.Range("M138").Value = -5782.7999
$ws.Range("N138").Value = -32938.2938

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 75111.14
$ws.Range("I2").Value = 128326
$ws.Range("J2").Value = 4158
$ws.Range("K2").Value = 128326
$ws.Range("L2").Value = 4158
$ws.Range("M2").Value = -128213
$ws.Range("N2").Value = -4384
# Row 61
$ws.Range("H61").Value = 3790.76
$ws.Range("I61").Value = 3206.25
$ws.Range("K61").Value = 3206.25
$ws.Range("M61").Value = -2994.25
# Row 110
$ws.Range("H110").Value = 267340.47
$ws.Range("I110").Value = 359973.16
$ws.Range("K110").Value = 359973.16
$ws.Range("M110").Value = -357928.16
# Row 116
$ws.Range("H116").Value = 75111.14
$ws.Range("I116").Value = 128326
$ws.Range("J116").Value = 4158
$ws.Range("K116").Value = 128326
$ws.Range("L116").Value = 4158
$ws.Range("M116").Value = -126032
$ws.Range("N116").Value = -8746
# Row 136
$ws.Range("H136").Value = 3790.76
$ws.Range("I136").Value = 3206.25
$ws.Range("K136").Value = 9618.75
$ws.Range("M136").Value = -7068.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 75111.14
$ws.Range("I3").Value = 128326
$ws.Range("J3").Value = 4158
$ws.Range("K3").Value = 128326
$ws.Range("L3").Value = 4158
$ws.Range("M3").Value = -128212
$ws.Range("N3").Value = -4386
# Row 54
$ws.Range("H54").Value = 13689
$ws.Range("I54").Value = 1585.3334
$ws.Range("J54").Value = 50000
$ws.Range("K54").Value = 1585.3334
$ws.Range("L54").Value = 50000
$ws.Range("M54").Value = -1101.3334
$ws.Range("N54").Value = -50968
# Row 97
$ws.Range("H97").Value = 10007
$ws.Range("I97").Value = 10007
$ws.Range("K97").Value = 10007
$ws.Range("M97").Value = -9016
# Row 122
$ws.Range("H122").Value = 82666.664
$ws.Range("J122").Value = 82666.664
$ws.Range("L122").Value = 82666.664
$ws.Range("N122").Value = -92466.664
# Row 134
$ws.Range("H134").Value = 36880.6
$ws.Range("I134").Value = 2360.5881
$ws.Range("K134").Value = 7081.7643
$ws.Range("M134").Value = -4546.7643

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 26
$ws.Range("H26").Value = 4000
$ws.Range("J26").Value = 4000
$ws.Range("L26").Value = 4000
$ws.Range("N26").Value = -4574
# Row 58
$ws.Range("H58").Value = 4317.1055
$ws.Range("I58").Value = 2184.182
$ws.Range("K58").Value = 2184.182
$ws.Range("M58").Value = -1981.182
# Row 62
$ws.Range("H62").Value = 2930.9167
$ws.Range("I62").Value = 2682.1428
$ws.Range("K62").Value = 2682.1428
$ws.Range("M62").Value = -2058.1428
# Row 65
$ws.Range("H65").Value = 2930.9167
$ws.Range("I65").Value = 2682.1428
$ws.Range("K65").Value = 13410.714
$ws.Range("M65").Value = -10290.714
# Row 74
$ws.Range("H74").Value = 79888
$ws.Range("J74").Value = 79888
$ws.Range("L74").Value = 79888
$ws.Range("N74").Value = -81636
# Row 77
$ws.Range("H77").Value = 79888
$ws.Range("J77").Value = 79888
$ws.Range("L77").Value = 239664
$ws.Range("N77").Value = -248400
# Row 132
$ws.Range("H132").Value = 2952.16
$ws.Range("I132").Value = 2306.1
$ws.Range("K132").Value = 6918.299999999999
$ws.Range("M132").Value = -4388.299999999999
# Row 134
$ws.Range("H134").Value = 288760.97
$ws.Range("I134").Value = 2621.4
$ws.Range("J134").Value = 1004109.9
$ws.Range("K134").Value = 7864.200000000001
$ws.Range("L134").Value = 3012329.7
$ws.Range("M134").Value = -5329.200000000001
$ws.Range("N134").Value = -3017399.7
# Row 136
$ws.Range("H136").Value = 4317.1055
$ws.Range("I136").Value = 2184.182
$ws.Range("K136").Value = 6552.545999999999
$ws.Range("M136").Value = -4002.545999999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 119995
$ws.Range("J37").Value = 119995
$ws.Range("L37").Value = 359985
$ws.Range("N37").Value = -360209
# Row 68
$ws.Range("H68").Value = 2090.7551
$ws.Range("J68").Value = 2081.3635
$ws.Range("L68").Value = 6244.0905
$ws.Range("N68").Value = -7866.0905
# Row 71
$ws.Range("H71").Value = 2090.7551
$ws.Range("J71").Value = 2081.3635
$ws.Range("L71").Value = 18732.2715
$ws.Range("N71").Value = -26844.2715
# Row 121
$ws.Range("H121").Value = 628009.1
$ws.Range("J121").Value = 912858.8
$ws.Range("L121").Value = 2738576.4
$ws.Range("N121").Value = -2741196.4
# Row 123
$ws.Range("H123").Value = 12633.333
$ws.Range("J123").Value = 12633.333
$ws.Range("L123").Value = 37899.999
$ws.Range("N123").Value = -42799.999
# Row 131
$ws.Range("H131").Value = 90794.69500000001
$ws.Range("J131").Value = 76852.21000000001
$ws.Range("L131").Value = 230556.63
$ws.Range("N131").Value = -240636.63

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 568209.2
$ws.Range("J80").Value = 636367.75
$ws.Range("L80").Value = 636367.75
$ws.Range("N80").Value = -638363.75
# Row 83
$ws.Range("H83").Value = 568209.2
$ws.Range("J83").Value = 636367.75
$ws.Range("L83").Value = 3181838.75
$ws.Range("N83").Value = -3191822.75
# Row 113
$ws.Range("H113").Value = 40417416
$ws.Range("I113").Value = 649992.9399999999
$ws.Range("K113").Value = 649992.9399999999
$ws.Range("M113").Value = -647822.9399999999
# Row 132
$ws.Range("H132").Value = 33031.03
$ws.Range("I132").Value = 2495.1365
$ws.Range("K132").Value = 7485.4095
$ws.Range("M132").Value = -4955.4095

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 291151.12
$ws.Range("I7").Value = 504390.9
$ws.Range("K7").Value = 504390.9
$ws.Range("M7").Value = -504278.9
# Row 61
$ws.Range("H61").Value = 4826.853
$ws.Range("I61").Value = 3744.52
$ws.Range("K61").Value = 3744.52
$ws.Range("M61").Value = -3542.52
# Row 68
$ws.Range("H68").Value = 59250.332
$ws.Range("I68").Value = 3619
$ws.Range("K68").Value = 3619
$ws.Range("M68").Value = -2870
# Row 71
$ws.Range("H71").Value = 59250.332
$ws.Range("I71").Value = 3619
$ws.Range("K71").Value = 18095
$ws.Range("M71").Value = -14351
# Row 82
$ws.Range("H82").Value = 2935.5715
$ws.Range("I82").Value = 1430.5
$ws.Range("K82").Value = 1430.5
$ws.Range("M82").Value = -1069.5
# Row 85
$ws.Range("H85").Value = 2935.5715
$ws.Range("I85").Value = 1430.5
$ws.Range("K85").Value = 1430.5
$ws.Range("M85").Value = -182.5
# Row 113
$ws.Range("H113").Value = 4826.853
$ws.Range("I113").Value = 3744.52
$ws.Range("K113").Value = 3744.52
$ws.Range("M113").Value = -1574.52
# Row 122
$ws.Range("H122").Value = 329208.62
$ws.Range("I122").Value = 5643.6665
$ws.Range("J122").Value = 693219.2
$ws.Range("K122").Value = 16930.9995
$ws.Range("L122").Value = 2079657.6
$ws.Range("M122").Value = -14480.9995
$ws.Range("N122").Value = -2084557.6
# Row 126
$ws.Range("H126").Value = 291151.12
$ws.Range("I126").Value = 504390.9
$ws.Range("K126").Value = 1513172.7
$ws.Range("M126").Value = -1510702.7

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 7181.727
$ws.Range("I62").Value = 6599.8
$ws.Range("K62").Value = 6599.8
$ws.Range("M62").Value = -5975.8
# Row 65
$ws.Range("H65").Value = 7181.727
$ws.Range("I65").Value = 6599.8
$ws.Range("K65").Value = 32999
$ws.Range("M65").Value = -29879
# Row 100
$ws.Range("H100").Value = 1168.5
$ws.Range("I100").Value = 1600
$ws.Range("K100").Value = 3200
$ws.Range("M100").Value = -2659
# Row 132
$ws.Range("H132").Value = 18530.287
$ws.Range("I132").Value = 2873.9185
$ws.Range("J132").Value = 63657.47
$ws.Range("K132").Value = 8621.755500000001
$ws.Range("L132").Value = 190972.41
$ws.Range("M132").Value = -6091.755500000001
$ws.Range("N132").Value = -196032.41
# Row 136
$ws.Range("H136").Value = 467007.12
$ws.Range("I136").Value = 670447.4
$ws.Range("K136").Value = 2011342.2
$ws.Range("M136").Value = -2008792.2
